# Updated cryptos list on Tue Jul 25 13:45:45 UTC 2023 with GitHub Actions
#
# This mirrors the daily cryptos-list refresh: the Price (column D) and
# Volume(1h) (column E) columns are refreshed with the latest scraped
# values. A couple of coins (ShibaInu / BitcoinCash) also swapped rank
# positions (rows 19/20) as part of the refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($Cell, $Text) {
    # Force the cell to be treated as text so that values which look
    # numeric (e.g. "0.9980", "238.08") keep their exact original
    # formatting instead of being normalized/rounded by Excel.
    $r = $ws.Range($Cell)
    $r.NumberFormat = "@"
    $r.Value = $Text
}

# Row 2 - Bitcoin
Set-PriceText "D2" "29.247.01"
$ws.Range("E2").Value = "  +0.08%  "

# Row 3 - Ethereum
Set-PriceText "D3" "1.861.82"
$ws.Range("E3").Value = "  +0.53%  "

# Row 4 - TetherUSD
Set-PriceText "D4" "0.9980"
$ws.Range("E4").Value = "  -0.26%  "

# Row 5 - XRP
Set-PriceText "D5" "0.6988"
$ws.Range("E5").Value = "  +0.62%  "

# Row 6 - BNB
Set-PriceText "D6" "236.53"
$ws.Range("E6").Value = "  -0.79%  "

# Row 7 - USDC
Set-PriceText "D7" "0.9985"
$ws.Range("E7").Value = "  -0.19%  "

# Row 8 - Dogecoin
Set-PriceText "D8" "0.07671"
$ws.Range("E8").Value = "  +1.19%  "

# Row 9 - Cardano
Set-PriceText "D9" "0.3041"
$ws.Range("E9").Value = "  -0.75%  "

# Row 10 - Solana
Set-PriceText "D10" "23.27"
$ws.Range("E10").Value = "  -1.07%  "

# Row 11 - TRON
Set-PriceText "D11" "0.08114"
$ws.Range("E11").Value = "  +0.34%  "

# Row 12 - WrappedEther
Set-PriceText "D12" "1.919.25"
$ws.Range("E12").Value = "  +3.68%  "

# Row 13 - Polygon
Set-PriceText "D13" "0.7160"
$ws.Range("E13").Value = "  -0.96%  "

# Row 14 - Polkadot (only Price changes)
Set-PriceText "D14" "5.150"

# Row 15 - Litecoin
Set-PriceText "D15" "89.64"
$ws.Range("E15").Value = "  +0.71%  "

# Row 16 - WrappedBTC
Set-PriceText "D16" "29.217.85"
$ws.Range("E16").Value = "  +0.00%  "

# Row 17 - Uniswap
Set-PriceText "D17" "5.738"
$ws.Range("E17").Value = "  -0.53%  "

# Row 18 - Avalanche (only Price changes)
Set-PriceText "D18" "13.17"

# Rows 19 & 20 - ShibaInu and BitcoinCash swap positions, with refreshed data
Set-PriceText "D19" "238.08"
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("E19").Value = "  -1.22%  "

Set-PriceText "D20" "0.000007720"
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("E20").Value = "  -0.01%  "

# Row 21 - Dai
Set-PriceText "D21" "0.9992"
$ws.Range("E21").Value = "  -0.02%  "

# Row 22 - WrappedliquidstakedEther2.0
Set-PriceText "D22" "2.108.23"
$ws.Range("E22").Value = "  +0.69%  "

# Row 23 - BinanceUSD
Set-PriceText "D23" "0.9985"
$ws.Range("E23").Value = "  -0.25%  "

# Row 24 - Chainlink
Set-PriceText "D24" "7.466"
$ws.Range("E24").Value = "  -1.83%  "

# Row 25 - Cosmos
Set-PriceText "D25" "9.018"
$ws.Range("E25").Value = "  -0.01%  "

# Row 26 - Monero
Set-PriceText "D26" "161.72"
$ws.Range("E26").Value = "  +0.06%  "

# Row 27 - Stellar
Set-PriceText "D27" "0.1451"
$ws.Range("E27").Value = "  -0.13%  "

# Row 28 - EthereumClassic
Set-PriceText "D28" "18.05"
$ws.Range("E28").Value = "  +0.09%  "

# Row 29 - LidoDAOToken
Set-PriceText "D29" "1.974"
$ws.Range("E29").Value = "  +2.11%  "

# Row 30 - Toncoin
Set-PriceText "D30" "1.408"
$ws.Range("E30").Value = "  +0.69%  "

# Row 31 - Filecoin
Set-PriceText "D31" "4.450"
$ws.Range("E31").Value = "  +0.49%  "

# Row 32 - PancakeSwap
Set-PriceText "D32" "1.479"
$ws.Range("E32").Value = "  -1.52%  "

# Row 33 - InternetComputer(DFINITY)
Set-PriceText "D33" "3.993"
$ws.Range("E33").Value = "  -1.32%  "

# Row 34 - Hedera
Set-PriceText "D34" "0.05189"
$ws.Range("E34").Value = "  -0.62%  "

# Row 35 - ARBITRUM
Set-PriceText "D35" "1.170"
$ws.Range("E35").Value = "  -1.81%  "

# Row 36 - ImmutableX
Set-PriceText "D36" "0.7081"
$ws.Range("E36").Value = "  +0.15%  "

# Row 37 - Frax
Set-PriceText "D37" "0.9973"
$ws.Range("E37").Value = "  -0.37%  "

# Row 38 - HuobiToken
Set-PriceText "D38" "2.647"
$ws.Range("E38").Value = "  -0.70%  "

# Row 39 - VeChain
Set-PriceText "D39" "0.01850"
$ws.Range("E39").Value = "  -0.57%  "

# Row 40 - MXToken (only Volume changes)
$ws.Range("E40").Value = "  +0.66%  "

# Row 41 - TrustWalletToken
Set-PriceText "D41" "0.9380"
$ws.Range("E41").Value = "  +2.30%  "

# Row 42 - Maker
Set-PriceText "D42" "1.133.42"
$ws.Range("E42").Value = "  +8.45%  "

# Row 43 - TheSandbox
Set-PriceText "D43" "0.4273"
$ws.Range("E43").Value = "  -0.42%  "

# Row 44 - Aave
Set-PriceText "D44" "70.97"
$ws.Range("E44").Value = "  +2.14%  "

# Row 45 - FraxShare
Set-PriceText "D45" "5.870"
$ws.Range("E45").Value = "  -1.53%  "

# Row 46 - PaxDollar
Set-PriceText "D46" "0.9987"
$ws.Range("E46").Value = "  -0.19%  "

# Row 47 - Quant
Set-PriceText "D47" "102.87"
$ws.Range("E47").Value = "  +0.52%  "

# Row 48 - RenderToken
Set-PriceText "D48" "1.801"
$ws.Range("E48").Value = "  +3.62%  "

# Row 49 - RocketPoolETH
Set-PriceText "D49" "2.007.11"
$ws.Range("E49").Value = "  +0.72%  "

# Row 50 - EnergySwap
Set-PriceText "D50" "9.152"
$ws.Range("E50").Value = "  -1.41%  "

# Row 51 - Aptos
Set-PriceText "D51" "6.954"
$ws.Range("E51").Value = "  -3.91%  "
